$wb = $excel.ActiveWorkbook

# --- Sheet: Overview ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-31 20:33:41"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-31 20:33:41"

# --- Sheet: zh-cn ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-08-31 20:33:37"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/98cda078-543e-48d5-b448-e610ef084672.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2be46dd678b30aa74b59bf004e6d8154988fda1c/e2e/98cda078-543e-48d5-b448-e610ef084672.md."
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-08-31 20:33:37"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2be46dd678b30aa74b59bf004e6d8154988fda1c/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- Sheet: de-de ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-08-31 20:33:41"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/98cda078-543e-48d5-b448-e610ef084672.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2be46dd678b30aa74b59bf004e6d8154988fda1c/e2e/98cda078-543e-48d5-b448-e610ef084672.md."
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-08-31 20:33:41"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2be46dd678b30aa74b59bf004e6d8154988fda1c/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
